$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the "Experimental" row's Value cell (B7) to the literal text "true".
# A direct Range.Value/Value2 assignment of "true" is auto-coerced to a
# native boolean cell (matching real Excel "typed in" semantics), which is
# not what we want here - the source data is a literal text string "true".
# Route it through a helper formula cell (whose computed result is the
# text string "true") and PasteSpecial just the values onto B7; a values-only
# paste of a formula's text result is not re-coerced to boolean, and the
# destination keeps its existing cell style.
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").ClearContents()

# Update the "Date" row's Value cell (B8) to the new timestamp
$ws.Range("B8").Value2 = "2025-01-28T15:58:19+00:00"
